$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# World Bank country classification refresh: FY2024 thresholds -> FY2023
# (latest classification compatible with other macroeconomic variables)

# year column: 2024 -> 2023 for all data rows
$ws.Range("F2").Value = 2023
$ws.Range("F3").Value = 2023
$ws.Range("F4").Value = 2023
$ws.Range("F5").Value = 2023

# updated income thresholds (min/max) for FY2023 classification
$ws.Range("D3").Value = 4465
$ws.Range("C4").Value = 4466
$ws.Range("D4").Value = 13845
$ws.Range("C5").Value = 13846

# mark the updated year cells (F3:F5) with an underline, as done by hand
$ws.Range("F3:F5").Font.Underline = $true

# leave the selection where the editor left it
$ws.Range("F3:F5").Select()
